# Improved the volatility point for better rating.
# Recomputed a handful of fundamental % inputs (and the two trend-score
# floating point values) for NXPI / AXISCADES.NS / CARTRADE.NS / CHOLAFIN.NS,
# which ripples into the Final Rating and the "Calculation Details" narration
# on both worksheets.

$wb = $excel.ActiveWorkbook
$wsRankings = $wb.Worksheets.Item("Rankings")
$wsDetails  = $wb.Worksheets.Item("Calculation Details")

# ---------------------------------------------------------------------------
# 1) Simple numeric / short-text cell updates on the "Rankings" sheet
# ---------------------------------------------------------------------------

# Row 2 - NXPI (only the two trend/RS scores move)
$wsRankings.Range("L2").Value2 = 35.45860667648925
$wsRankings.Range("Q2").Value2 = 30.38357733960528

# Row 3 - AXISCADES.NS
$wsRankings.Range("B3").Value2 = 24.3
$wsRankings.Range("C3").Value2 = "59.4% -> 100"
$wsRankings.Range("D3").Value2 = "3.0% -> 12.0"
$wsRankings.Range("G3").Value2 = "34.0% -> 100"
$wsRankings.Range("H3").Value2 = "85.0% -> 100"

# Row 4 - CARTRADE.NS
$wsRankings.Range("B4").Value2 = 29.09
$wsRankings.Range("C4").Value2 = "61.0% -> 100"
$wsRankings.Range("D4").Value2 = "39.0% -> 100"
$wsRankings.Range("G4").Value2 = "78.0% -> 100"
$wsRankings.Range("H4").Value2 = "35.0% -> 100"

# Row 5 - CHOLAFIN.NS
$wsRankings.Range("B5").Value2 = 40.53
$wsRankings.Range("C5").Value2 = "49.93% -> 99.86"
$wsRankings.Range("D5").Value2 = "42.0% -> 100"
$wsRankings.Range("G5").Value2 = "27.0% -> 100"
$wsRankings.Range("H5").Value2 = "23.0% -> 100"

# ---------------------------------------------------------------------------
# 2) Per-ticker substring replacements that ripple through the long
#    "Calculation Details" narration cells (column R on Rankings, column A
#    on the Calculation Details sheet). Every occurrence of each "metric:
#    old -> Score: old" phrase is swapped for the corresponding new phrase.
# ---------------------------------------------------------------------------

function Apply-Replacements {
    param($text, $pairs)
    $result = $text
    foreach ($pair in $pairs) {
        $result = $result.Replace($pair[0], $pair[1])
    }
    return $result
}

$nxpiPairs = @(
    ,@('35.45860667649082', '35.45860667648925')
    ,@('30.383570166515007', '30.38357733960528')
)

$axiscadesPairs = @(
    ,@('Promoter Holding: 50.6% -> Score: 100', 'Promoter Holding: 59.4% -> Score: 100')
    ,@('Institutional Holding: 8.7% -> Score: 34.8', 'Institutional Holding: 3.0% -> Score: 12.0')
    ,@('Profit Growth YoY: 38.5% -> Score: 100', 'Profit Growth YoY: 34.0% -> Score: 100')
    ,@('Profit CAGR 5Y: 11.5% -> Score: 57.5', 'Profit CAGR 5Y: 85.0% -> Score: 100')
    ,@('Final Rating: 23.69', 'Final Rating: 24.3')
)

$cartradePairs = @(
    ,@('Promoter Holding: 60.0% -> Score: 100', 'Promoter Holding: 61.0% -> Score: 100')
    ,@('Institutional Holding: 10.0% -> Score: 40.0', 'Institutional Holding: 39.0% -> Score: 100')
    ,@('Profit Growth YoY: 20.0% -> Score: 100', 'Profit Growth YoY: 78.0% -> Score: 100')
    ,@('Profit CAGR 5Y: 25.0% -> Score: 100', 'Profit CAGR 5Y: 35.0% -> Score: 100')
    ,@('Final Rating: 27.07', 'Final Rating: 29.09')
)

$cholafinPairs = @(
    ,@('Promoter Holding: 50.0% -> Score: 100', 'Promoter Holding: 49.93% -> Score: 99.86')
    ,@('Institutional Holding: 15.0% -> Score: 60.0', 'Institutional Holding: 42.0% -> Score: 100')
    ,@('Profit Growth YoY: 10.0% -> Score: 50.0', 'Profit Growth YoY: 27.0% -> Score: 100')
    ,@('Profit CAGR 5Y: 12.0% -> Score: 60.0', 'Profit CAGR 5Y: 23.0% -> Score: 100')
    ,@('Final Rating: 34.24', 'Final Rating: 40.53')
)

# --- Rankings!R2:R5 (one long narration cell per ticker row) --------------
$wsRankings.Range("R2").Value2 = Apply-Replacements $wsRankings.Range("R2").Value2 $nxpiPairs
$wsRankings.Range("R3").Value2 = Apply-Replacements $wsRankings.Range("R3").Value2 $axiscadesPairs
$wsRankings.Range("R4").Value2 = Apply-Replacements $wsRankings.Range("R4").Value2 $cartradePairs
$wsRankings.Range("R5").Value2 = Apply-Replacements $wsRankings.Range("R5").Value2 $cholafinPairs

# --- Calculation Details!A2:A681 (one narration line per cell) -----------
$usedRange = $wsDetails.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $wsDetails.Cells.Item($r, 1)
    $val = $cell.Value2
    if ($null -eq $val) { continue }

    if ($val.StartsWith("NXPI:")) {
        $newVal = Apply-Replacements $val $nxpiPairs
    } elseif ($val.StartsWith("AXISCADES.NS:")) {
        $newVal = Apply-Replacements $val $axiscadesPairs
    } elseif ($val.StartsWith("CARTRADE.NS:")) {
        $newVal = Apply-Replacements $val $cartradePairs
    } elseif ($val.StartsWith("CHOLAFIN.NS:")) {
        $newVal = Apply-Replacements $val $cholafinPairs
    } else {
        $newVal = $val
    }

    if ($newVal -ne $val) {
        $cell.Value2 = $newVal
    }
}
